# "change names and units, by Anita"
# - rename a few entries in the "name" column (C) to be more descriptive
# - fix the "unit" for Soil Organic Matter (D10: "???" -> "%")
# - add a new "measurementMethod" column (E) describing how each variable
#   was measured, including its header in E1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C (name) renames ---
$ws.Range("C2").Value  = "pH measured in KCl"      # was "pH KCl"
$ws.Range("C5").Value  = "Sulphate content"        # was "Sulphate"
$ws.Range("C6").Value  = "Bicarbonate content"     # was "Bicarbonate"
$ws.Range("C7").Value  = "Content of total N"      # was "Total N"

# --- Column D (unit) fix ---
$ws.Range("D10").Value = "%"                       # was "???"

# --- New column E: measurementMethod ---
$ws.Range("E1").Value  = "measurementMethod"
$ws.Range("E2").Value  = "Electrochemical method"
$ws.Range("E3").Value  = "Electrochemical method"
$ws.Range("E4").Value  = "Electrochemical method"
$ws.Range("E5").Value  = "Ion chromatography"
$ws.Range("E6").Value  = "Volumetric method"
$ws.Range("E7").Value  = "Dumas method"
$ws.Range("E8").Value  = "ICP-MS"
$ws.Range("E9").Value  = "Spectrophotometry"
$ws.Range("E10").Value = "Volumetric method"
$ws.Range("E11").Value = "Calculation"
$ws.Range("E12").Value = "Densimetry"
$ws.Range("E13").Value = "Densimetry"
$ws.Range("E14").Value = "Densimetry"
$ws.Range("E15").Value = "Gravimetry"

# Widen the new column to fit its contents, like the other labelled columns
$ws.Range("E1").EntireColumn.ColumnWidth = 24.14

# Reflect the view state captured in the saved workbook: zoomed out a bit
# and the active selection left on C7.
$excel.ActiveWindow.Zoom = 176
[void]$ws.Range("C7").Select()
